# Apply content updates to the LOB1254 course-description workbook.
# Columns B and C on each row hold the same text (current vs. modified copy),
# so both need to be updated together to keep them in sync with the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-BothCols {
    param($Row, $Text)
    $ws.Cells.Item($Row, 2).Value = $Text
    $ws.Cells.Item($Row, 3).Value = $Text
}

# Name: "Geology for Environmental Engineering" -> "Geology"
Set-BothCols 4 "Geology"

# Ativação: "01/01/2020" -> "01/01/2022"
# Force text format first so the date-like string isn't auto-converted
# into a date serial number (the source data stores it as plain text).
$ws.Cells.Item(8, 2).NumberFormat = "@"
$ws.Cells.Item(8, 3).NumberFormat = "@"
Set-BothCols 8 "01/01/2022"

# Objectives:
Set-BothCols 11 "Provide basic knowledge about terrestrial materials and the main geological processes."

# Programa resumido:
Set-BothCols 14 "Processos endógenos e exógenos da Terra. Materiais constituintes da crosta terrestre (minerais e rochas)."

# Short syllabus:
Set-BothCols 15 "Endogenous and exogenous processes of the Earth. Materials constituting the earth's crust (minerals and rocks)."

# Programa:
Set-BothCols 16 "Breve história da Geologia. Materiais constituintes da crosta terrestre (minerais e rochas). Origem e constituição do universo, do sistema solar e da Terra. Estrutura interna da Terra. Composição da Terra. Processos endógenos e exógenos (dinâmica interna e externa da Terra).  Teoria da tectônica de placas.  Rochas ígneas e vulcanismo. Rochas metamórficas e metamorfismo. Rochas sedimentares. Intemperismo, erosão, transporte de sedimentos.  Estrutura geológicas. Tempo geológico e estratigrafia."

# Syllabus:
Set-BothCols 17 "Brief history of geology. Materials constituting the earth's crust (minerals and rocks). Origin and constitution of the universe, the solar system and the earth. Internal structure of the earth. Composition of the earth. Endogenous and exogenous processes (internal and external dynamics of the earth). Plate tectonics theory.  Igneous rocks and vulcanismo. Metamorphic rocks and metamorphism. Sedimentary rocks. Weathering, erosion, sediment transport. Geological structure. Geological time and stratigraphy."

# Método:
Set-BothCols 19 "Aulas teóricas expositivas, atividades individuais e em grupo, relatórios e provas."

# Critério:
Set-BothCols 20 "Média ponderada de provas  e atividades."

# Norma de recuperação:
Set-BothCols 21 "1 (uma) prova escrita"

# Bibliografia:
Set-BothCols 22 "Bibliografia básica:PRESS, F.; SIEVER, R.; GROTZINGER, J.; JORDAN, T. H. Para entender a Terra. Porto Alegre: Bookman, 2008. 656p.REED, W.; MONROE, J. S. Fundamentos de Geologia. São Paulo: Cengage Learning, 2011. 508p.Bibliografia complementar:TEIXEIRA, W.; FAIRCHILD, T. R.; DE TOLEDO, M. C. M.; TAIOLI, F. Decifrando a Terra. São Paulo: Companhia Editora Nacional, 2003. 623p."
